$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.45
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 5.1
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 8.6
$ws.Range("N2").Value = 1.91
$ws.Range("O2").Value = 1.23
$ws.Range("P2").Value = 1.91
$ws.Range("Q2").Value = 1.61
$ws.Range("R2").Value = 1.35
$ws.Range("S2").Value = 2.52
$ws.Range("V2").Value = 1.04
$ws.Range("W2").Value = 1.11
$ws.Range("F3").Value = 1.74
$ws.Range("G3").Value = 1.86
$ws.Range("H3").Value = 5.4
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 3.45
$ws.Range("K3").Value = 3.75
$ws.Range("R3").Value = 1.26
$ws.Range("T3").Value = 2.06
$ws.Range("U3").Value = 1.79
$ws.Range("AF3").Value = 11
$ws.Range("AG3").Value = 10.5
$ws.Range("F4").Value = 1.4
$ws.Range("G4").Value = 1.41
$ws.Range("H4").Value = 8.6
$ws.Range("I4").Value = 9.4
$ws.Range("J4").Value = 5.6
$ws.Range("K4").Value = 5.9
$ws.Range("L4").Value = 1.29
$ws.Range("N4").Value = 5.9
$ws.Range("P4").Value = 2.68
$ws.Range("Q4").Value = 1.49
$ws.Range("R4").Value = 1.68
$ws.Range("S4").Value = 2.22
$ws.Range("U4").Value = 2.08
$ws.Range("V4").Value = 1.12
$ws.Range("W4").Value = 3.35
$ws.Range("X4").Value = 27
$ws.Range("Y4").Value = 40
$ws.Range("AB4").Value = 12
$ws.Range("AE4").Value = 120
$ws.Range("AH4").Value = 24
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 13
$ws.Range("AK4").Value = 13.5
$ws.Range("AM4").Value = 110
$ws.Range("AN4").Value = 4.7
$ws.Range("AO4").Value = 120
$ws.Range("F5").Value = 4
$ws.Range("H5").Value = 1.93
$ws.Range("I5").Value = 2.04
$ws.Range("J5").Value = 3.7
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.36
$ws.Range("N5").Value = 3.9
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 1.79
$ws.Range("R5").Value = 1.4
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 1.72
$ws.Range("U5").Value = 2.14
$ws.Range("V5").Value = 1.96
$ws.Range("AB5").Value = 20
$ws.Range("F6").Value = 6.8
$ws.Range("G6").Value = 7.4
$ws.Range("H6").Value = 1.46
$ws.Range("I6").Value = 1.48
$ws.Range("J6").Value = 5.3
$ws.Range("K6").Value = 5.6
$ws.Range("L6").Value = 1.2
$ws.Range("M6").Value = 1.03
$ws.Range("Q6").Value = 1.44
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 1.67
$ws.Range("U6").Value = 2.26
$ws.Range("V6").Value = 3.05
$ws.Range("W6").Value = 1.16
$ws.Range("Z6").Value = 14.5
$ws.Range("AA6").Value = 14.5
$ws.Range("AF6").Value = 70
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 65
$ws.Range("AO6").Value = 4.9
$ws.Range("F7").Value = 3.95
$ws.Range("G7").Value = 5.3
$ws.Range("H7").Value = 1.83
$ws.Range("I7").Value = 2.08
$ws.Range("J7").Value = 3.4
$ws.Range("K7").Value = 4.2
$ws.Range("L7").Value = 1.45
$ws.Range("N7").Value = 3.2
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 1.75
$ws.Range("Q7").Value = 2
$ws.Range("S7").Value = 4.1
$ws.Range("T7").Value = 1.88
$ws.Range("V7").Value = 1.92
$ws.Range("Y7").Value = 10
$ws.Range("AC7").Value = 10
$ws.Range("G8").Value = 2.72
$ws.Range("H8").Value = 2.72
$ws.Range("I8").Value = 3.05
$ws.Range("L8").Value = 1.31
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 4.5
$ws.Range("Q8").Value = 1.65
$ws.Range("S8").Value = 2.62
$ws.Range("V8").Value = 1.49
$ws.Range("W8").Value = 1.59
$ws.Range("F9").Value = 5.3
$ws.Range("G9").Value = 7.2
$ws.Range("H9").Value = 1.64
$ws.Range("I9").Value = 1.76
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4.5
$ws.Range("L9").Value = 1.35
$ws.Range("M9").Value = 1.05
$ws.Range("S9").Value = 2.98
$ws.Range("V9").Value = 2.32
$ws.Range("F10").Value = 1.45
$ws.Range("G10").Value = 1.51
$ws.Range("I10").Value = 10
$ws.Range("J10").Value = 4.4
$ws.Range("L10").Value = 1.35
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 3.9
$ws.Range("O10").Value = 1.27
$ws.Range("Q10").Value = 1.79
$ws.Range("S10").Value = 3.05
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 1.8
$ws.Range("W10").Value = 2.92
$ws.Range("F11").Value = 1.7
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 4.5
$ws.Range("I11").Value = 5.2
$ws.Range("J11").Value = 4.5
$ws.Range("O11").Value = 1.17
$ws.Range("P11").Value = 2.5
$ws.Range("R11").Value = 1.62
$ws.Range("V11").Value = 1.24
$ws.Range("W11").Value = 2.32
$ws.Range("Y11").Value = 28
$ws.Range("Z11").Value = 44
$ws.Range("AB11").Value = 14
$ws.Range("AD11").Value = 19.5
$ws.Range("AF11").Value = 13.5
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 19.5
$ws.Range("AK11").Value = 16
$ws.Range("AM11").Value = 70
$ws.Range("AN11").Value = 7.2
$ws.Range("AO11").Value = 38
$ws.Range("G12").Value = 1.73
$ws.Range("H12").Value = 6.4
$ws.Range("I12").Value = 8.199999999999999
$ws.Range("J12").Value = 3.35
$ws.Range("K12").Value = 4.2
$ws.Range("L12").Value = 1.52
$ws.Range("N12").Value = 2.68
$ws.Range("S12").Value = 4.9
$ws.Range("T12").Value = 2.28
$ws.Range("W12").Value = 2.36
$ws.Range("AB12").Value = 1000
$ws.Range("F13").Value = 3.1
$ws.Range("I13").Value = 2.82
$ws.Range("J13").Value = 3.1
$ws.Range("L13").Value = 1.55
$ws.Range("M13").Value = 1.11
$ws.Range("N13").Value = 2.64
$ws.Range("O13").Value = 1.53
$ws.Range("Q13").Value = 2.56
$ws.Range("T13").Value = 2.04
$ws.Range("V13").Value = 1.55
$ws.Range("X13").Value = 9.4
$ws.Range("Y13").Value = 17
$ws.Range("Z13").Value = 16.5
$ws.Range("AA13").Value = 46
$ws.Range("AD13").Value = 13.5
$ws.Range("AE13").Value = 40
$ws.Range("AF13").Value = 19.5
$ws.Range("AG13").Value = 18
$ws.Range("AH13").Value = 25
$ws.Range("AI13").Value = 70
$ws.Range("AJ13").Value = 70
$ws.Range("AK13").Value = 55
$ws.Range("AM13").Value = 200
$ws.Range("AN13").Value = 70
$ws.Range("F14").Value = 1.91
$ws.Range("G14").Value = 1.97
$ws.Range("H14").Value = 4.7
$ws.Range("I14").Value = 5.2
$ws.Range("K14").Value = 3.6
$ws.Range("M14").Value = 1.1
$ws.Range("P14").Value = 1.67
$ws.Range("R14").Value = 1.24
$ws.Range("V14").Value = 1.24
$ws.Range("W14").Value = 2.02
$ws.Range("Y14").Value = 17.5
$ws.Range("Z14").Value = 34
$ws.Range("AA14").Value = 130
$ws.Range("AD14").Value = 21
$ws.Range("AE14").Value = 80
$ws.Range("AI14").Value = 980
$ws.Range("AK14").Value = 25
$ws.Range("AL14").Value = 50
$ws.Range("AM14").Value = 180
$ws.Range("AN14").Value = 23
$ws.Range("L15").Value = 1.56
$ws.Range("M15").Value = 1.11
$ws.Range("S15").Value = 5
$ws.Range("W15").Value = 1.72
$ws.Range("Y15").Value = 12.5
$ws.Range("Z15").Value = 27
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AD15").Value = 19.5
$ws.Range("AJ15").Value = 34
$ws.Range("AK15").Value = 32
$ws.Range("AN15").Value = 38
$ws.Range("AO15").Value = 100
$ws.Range("H16").Value = 3.4
$ws.Range("I16").Value = 3.8
$ws.Range("J16").Value = 3.05
$ws.Range("Q16").Value = 2.6
$ws.Range("S16").Value = 5.6
$ws.Range("U16").Value = 1.74
$ws.Range("V16").Value = 1.37
$ws.Range("Y16").Value = 10.5
$ws.Range("AF16").Value = 16.5
$ws.Range("AG16").Value = 13
$ws.Range("AJ16").Value = 1000
$ws.Range("AK16").Value = 36
$ws.Range("AM16").Value = 230
Write-Output "Applied 236 cell updates"
